$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; this shifts header (was row1, style s=1) to row2, data down by 1
$ws.Rows("1:1").Insert()

# At this point row2 A:K has style s=1 (carried from old row1). Copy that format to new row1 A:K only.
$ws.Range("A2:K2").Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate new row1 with sequential index values 0..10
for ($c = 1; $c -le 11; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 1
}

# Remove style/formatting from row2 (the header text row), matching plain data rows
$ws.Rows("2:2").ClearFormats()

# In new row2, J2 and K2 should be blank (they previously held 'thread_size' / 'material_surface')
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
